$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Entry")

$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("Z4").Value = 1

$fec = $wb.Worksheets.Item("Central Ontario FEC")
$fec.Range("Y4").Formula = "=G4+H4+I4+J4+K4+N4+'User Entry'!X4"
$fec.Range("N5").Formula = "=IF(AND(O4<30,P4<20,Q4<20,G4<20, F4<20,K4<20),N4^0.5,0)"

# View changes: clear the stale "Central Ontario FEC" selection (was W5) back to A1,
# then re-zoom + re-select "User Entry" (which stays the active/visible tab).
$fec.Activate()
$fec.Range("A1").Select()

$ws.Activate()
$excel.ActiveWindow.Zoom = 85
